$d = $word.ActiveDocument

$d.Content.Find.Execute("`${tax}", $true, $false, $false, $false, $false,
                         $true, 1, $false, "`${tax} %", 2)
